# Apply the benchmark-stats update to the single-column table.
# Each table row holds one value in its sole cell; we rewrite the cell
# text for the rows whose values changed. Rows 44-46 originally packed
# several tab-separated numbers into one run - they collapse down to
# just their first value.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($rowIndex, $text) {
    $cell = $t.Rows.Item($rowIndex).Cells.Item(1)
    $cell.Range.Text = $text
}

Set-CellText 1 "0M"
Set-CellText 2 "0M"
Set-CellText 3 "0M"
Set-CellText 4 "151"

Set-CellText 5 "0.00002"
Set-CellText 6 "0.00010"
Set-CellText 8 "0.00002"
Set-CellText 11 "0.00010"
Set-CellText 12 "0.00602"

Set-CellText 44 "99.99"
Set-CellText 45 "0.01"
Set-CellText 46 "100"
